$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ophold")

# Add the new "Opholdsregler" column to the Ophold table (extends A1:B25 -> A1:C25)
$table = $ws.ListObjects.Item("Ophold")
$table.ListColumns.Add() | Out-Null

# Header
$ws.Range("C1").Value = "Opholdsregler"

# Values: "EU/EØS" for rows whose Opholdsgrundlag (col A) is an EU/EØS category,
# "Ikke-EU/EØS" for everything else.
$euRows = @(11, 12, 13, 23)
for ($r = 2; $r -le 25; $r++) {
    if ($euRows -contains $r) {
        $ws.Range("C$r").Value = "EU/EØS"
    } else {
        $ws.Range("C$r").Value = "Ikke-EU/EØS"
    }
}

# Column width for the new column
$ws.Columns.Item(3).ColumnWidth = 14.15

# Sheet2 ("Ophold") becomes the active/selected sheet, with a new selection
$ws.Activate()
$ws.Range("F14").Select()
